$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# --- Reference cells for format copying (unchanged text/number/percent styles) ---
# C15 = style 13 (text placeholder), G15 = style 14 (count number), H15 = style 15 (percent number)

# Row 14
$ws.Range("D14").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 15
$ws.Range("N15").Value = 650
$ws.Range("H15").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16
$ws.Range("C16").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 55.555555555555
$ws.Range("I16").Value = 115
$ws.Range("J16").Value = 126
$ws.Range("K16").Value = -8.730158730158
$ws.Range("L16").Value = -15.441176470588
$ws.Range("M16").Value = -37.5
$ws.Range("N16").Value = 400

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 215
$ws.Range("J17").Value = 195
$ws.Range("K17").Value = 10.25641025641
$ws.Range("L17").Value = 5.392156862745
$ws.Range("M17").Value = 49.305555555555
$ws.Range("N17").Value = 1243.75

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 86
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -7.52688172043
$ws.Range("L18").Value = -5.494505494505
$ws.Range("M18").Value = 2.380952380952
$ws.Range("N18").Value = 514.285714285714

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 28.571428571428
$ws.Range("I19").Value = 294
$ws.Range("J19").Value = 283
$ws.Range("K19").Value = 3.886925795053
$ws.Range("L19").Value = 1.030927835051
$ws.Range("M19").Value = 41.346153846153
$ws.Range("N19").Value = 3575

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E20").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 12
$ws.Range("I20").Value = 116
$ws.Range("K20").Value = 30.337078651685
$ws.Range("L20").Value = -1.694915254237
$ws.Range("M20").Value = 107.142857142857
$ws.Range("N20").Value = 544.444444444444

# Row 21
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 10.958904109589
$ws.Range("I21").Value = 841
$ws.Range("J21").Value = 799
$ws.Range("K21").Value = 5.256570713391
$ws.Range("L21").Value = -1.522248243559
$ws.Range("M21").Value = 21.356421356421
$ws.Range("N21").Value = 938.271604938272

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("G15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F22").Value = 1
$ws.Range("G15").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = -58.333333333333
$ws.Range("L22").Value = -72.222222222222
$ws.Range("M22").Value = -61.538461538461

# Row 23
$ws.Range("D23").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("L23").Value = 15.384615384615

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 51
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = -37.037037037037
$ws.Range("I24").Value = 617
$ws.Range("J24").Value = 721
$ws.Range("K24").Value = -14.424410540915
$ws.Range("L24").Value = -17.623497997329
$ws.Range("M24").Value = 32.974137931034

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -65.217391304347
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 225
$ws.Range("K25").Value = -43.111111111111
$ws.Range("L25").Value = -26.857142857142

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 296
$ws.Range("J26").Value = 261
$ws.Range("K26").Value = 13.409961685823
$ws.Range("L26").Value = -2.631578947368
$ws.Range("M26").Value = -23.711340206185

# Row 27
$ws.Range("C27").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F27").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H27").Value = -100
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 14.285714285714
$ws.Range("L27").Value = 23.076923076923

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("G15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F28").Value = 3
$ws.Range("I28").Value = 44
$ws.Range("K28").Value = 57.142857142857
$ws.Range("L28").Value = 29.411764705882

# Row 29
$ws.Range("D29").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E29").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N29").Value = -20

# Row 30
$ws.Range("D30").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("N30").Value = -20
